# CardValue.xlsx - "added some weapon and card values"
# Adds three new columns (G: Init Level, H: Max Level, I: Damage) to the
# existing card table and appends a brand-new card row (row 16: AAA / AoE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for columns G:I (row 1) ---
$ws.Range("G1").Value = "Init Level"
$ws.Range("H1").Value = "Max Level"
$ws.Range("I1").Value = "Damage"

# --- Fill G:I for the existing data rows (2-15) ---
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 7).Value = 1   # G = Init Level
    $ws.Cells.Item($r, 8).Value = 3   # H = Max Level
    $ws.Cells.Item($r, 9).Value = 1   # I = Damage
}

# --- Brand-new row 16: AAA / AoE card ---
$ws.Range("A16").Formula = "=ROW()-2"
$ws.Range("B16").Value = "AAA"
$ws.Range("C16").Value = "AoE"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = "None"
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 1

# --- Column F is narrower now that it only holds short text ---
$ws.Columns.Item(6).ColumnWidth = 20.28

# --- Selection moved to the newly added I15:I16 cells ---
$ws.Range("I15:I16").Select()
